$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append the newest daily-tracker row (2025/10/04, Saturday) following the
# existing table's pattern: date + weekday kept as literal text, the two
# numeric columns as real numbers.
#
# Force text formatting on the date cell before assignment so "2025/10/04"
# is stored as a literal string instead of being auto-parsed into a date
# serial number (matching how the existing rows store their dates).
$ws.Range("A57").NumberFormat = "@"
$ws.Range("A57").Value = "2025/10/04"
$ws.Range("A57").Style = "Normal"

$ws.Range("B57").Value = "土"
$ws.Range("C57").Value = 4
$ws.Range("D57").Value = 38
